# Auto-generated Excel COM-interop edit script
# Updates the cryptos price list (Price column D, Volume(1h) column E,
# plus a name/link swap of rows 42-43) to match the target snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D receive numeric-looking text (e.g. "0.999", "45.35").
# Excel auto-converts such strings to numbers on assignment, so the
# destination cells are pre-formatted as Text ("@") before the write and
# the formatting is cleared again afterwards so no stray number format
# is left behind - only the cell value changes, exactly like the diff.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

# --- Price (column D) updates ---
$ws.Range("D2").Value = "67.716.38"
$ws.Range("D3").Value = "3.324.82"
$ws.Range("D5").Value = "581.05"
$ws.Range("D6").Value = "175.18"
$ws.Range("D7").Value = "0.999"
$ws.Range("D8").Value = "0.588"
$ws.Range("D9").Value = "3.321.91"
$ws.Range("D10").Value = "0.178"
$ws.Range("D11").Value = "0.576"
$ws.Range("D12").Value = "45.35"
$ws.Range("D13").Value = "0.0000269"
$ws.Range("D14").Value = "660.19"
$ws.Range("D15").Value = "3.868.09"
$ws.Range("D16").Value = "8.39"
$ws.Range("D17").Value = "67.631.21"
$ws.Range("D19").Value = "3.328.89"
$ws.Range("D20").Value = "17.34"
$ws.Range("D21").Value = "10.97"
$ws.Range("D22").Value = "0.886"
$ws.Range("D25").Value = "98.57"
$ws.Range("D26").Value = "3.85"
$ws.Range("D28").Value = "9.25"
$ws.Range("D29").Value = "33.34"
$ws.Range("D30").Value = "8.42"
$ws.Range("D31").Value = "7.19"
$ws.Range("D32").Value = "568.33"
$ws.Range("D33").Value = "10.92"
$ws.Range("D35").Value = "0.999"
$ws.Range("D36").Value = "56.57"
$ws.Range("D37").Value = "3.667.51"
$ws.Range("D38").Value = "3.29"
$ws.Range("D39").Value = "34.24"
$ws.Range("D41").Value = "2.62"
$ws.Range("D44").Value = "0.333"
$ws.Range("D45").Value = "0.0₃0661"
$ws.Range("D51").Value = "129.49"

# Row 42/43 Coin/Link/Price swap (ApeXProtocol <-> Stacks, refreshed values)
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Value = "3.11"
$ws.Range("B43").Value = "ApeXProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D43").Value = "3.36"

# Restore default (General) number formatting on column D now that all
# text values are committed.
$priceRange.ClearFormats()

# --- Volume(1h) (column E) updates ---
$ws.Range("E2").Value = "  -0.22%  "
$ws.Range("E3").Value = "  +0.05%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("E5").Value = "  -0.36%  "
$ws.Range("E6").Value = "  -4.18%  "
$ws.Range("E8").Value = "  -0.45%  "
$ws.Range("E9").Value = "  +0.15%  "
$ws.Range("E10").Value = "  -0.32%  "
$ws.Range("E11").Value = "  -0.58%  "
$ws.Range("E12").Value = "  -2.40%  "
$ws.Range("E13").Value = "  -1.63%  "
$ws.Range("E14").Value = "  +3.97%  "
$ws.Range("E15").Value = "  +0.26%  "
$ws.Range("E16").Value = "  -0.93%  "
$ws.Range("E18").Value = "  -0.80%  "
$ws.Range("E19").Value = "  +0.17%  "
$ws.Range("E20").Value = "  -2.08%  "
$ws.Range("E21").Value = "  +0.36%  "
$ws.Range("E22").Value = "  -1.84%  "
$ws.Range("E23").Value = "  +5.71%  "
$ws.Range("E24").Value = "  -3.74%  "
$ws.Range("E25").Value = "  +1.40%  "
$ws.Range("E26").Value = "  -3.79%  "
$ws.Range("E27").Value = "  -4.11%  "
$ws.Range("E28").Value = "  -3.51%  "
$ws.Range("E29").Value = "  +2.49%  "
$ws.Range("E30").Value = "  -2.37%  "
$ws.Range("E31").Value = "  +6.47%  "
$ws.Range("E32").Value = "  -4.27%  "
$ws.Range("E33").Value = "  -0.39%  "
$ws.Range("E34").Value = "  +0.06%  "
$ws.Range("E35").Value = "  +0.13%  "
$ws.Range("E36").Value = "  +1.59%  "
$ws.Range("E37").Value = "  -7.18%  "
$ws.Range("E38").Value = "  -6.68%  "
$ws.Range("E39").Value = "  +4.84%  "
$ws.Range("E40").Value = "  +0.56%  "
$ws.Range("E41").Value = "  -2.58%  "
$ws.Range("E44").Value = "  -1.77%  "
$ws.Range("E45").Value = "  -3.66%  "
$ws.Range("E46").Value = "  -1.98%  "
$ws.Range("E47").Value = "  +1.30%  "
$ws.Range("E48").Value = "  -0.85%  "
$ws.Range("E49").Value = "  -0.38%  "
$ws.Range("E50").Value = "  -1.85%  "
$ws.Range("E51").Value = "  -0.86%  "
$ws.Range("E42").Value = "  -4.50%  "
$ws.Range("E43").Value = "  -1.12%  "
